{"js": "// The document's single table tracks feature status per row: the first\n// column names the feature, and the \"Atteint\" column (index 1) holds an\n// \"x\" marker when that feature has been reached.\n//\n// This change marks five more rows as reached by writing a lowercase\n// \"x\" into their \"Atteint\" cell:\n//   - four rows whose \"Atteint\" cell was empty gain a new \"x\" run;\n//   - \"Tour \u00e0 4 joueurs\" already had an uppercase \"X\", normalized to \"x\".\n\nconst targetLabels = [\n  \"2 joueurs en r\u00e9seau\",\n  \"Affichage du plateau pour les 2 joueurs r\u00e9seau\",\n  \"Tour \u00e0 4 joueurs\",\n  \"4 joueurs en r\u00e9seau\",\n  \"Affichage du plateau pour les 4 joueurs r\u00e9seau\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Load the label (first column) text of every row so we can locate the\n// target rows by name rather than relying on a hard-coded row index.\nconst labelCells = [];\nfor (let i = 0; i < table.rowCount; i++) {\n  const cell = table.getCell(i, 0);\n  cell.body.load(\"text\");\n  labelCells.push(cell);\n}\nawait context.sync();\n\n// For each target row, grab the first paragraph of its \"Atteint\" cell\n// (column index 1) so we can inspect/update its text.\nconst targetParagraphSets = [];\nfor (let i = 0; i < table.rowCount; i++) {\n  const label = labelCells[i].body.text.trim();\n  if (targetLabels.indexOf(label) !== -1) {\n    const atteintCell = table.getCell(i, 1);\n    atteintCell.body.paragraphs.load(\"items\");\n    targetParagraphSets.push(atteintCell.body.paragraphs);\n  }\n}\nawait context.sync();\n\nfor (const paras of targetParagraphSets) {\n  paras.items[0].load(\"text\");\n}\nawait context.sync();\n\nfor (const paras of targetParagraphSets) {\n  const para = paras.items[0];\n  if (para.text && para.text.trim().length > 0) {\n    // Existing marker (e.g. \"X\") \u2014 normalize it to lowercase \"x\".\n    para.getRange().insertText(\"x\", Word.InsertLocation.replace);\n  } else {\n    // Empty cell paragraph \u2014 add the new \"x\" run, keeping the existing\n    // paragraph properties (center alignment) intact.\n    para.insertText(\"x\", Word.InsertLocation.end);\n  }\n}\nawait context.sync();\n", "ps1": "# The document's single table tracks feature status per row: the first\n# column names the feature, and the \"Atteint\" column (column 2) holds an\n# \"x\" marker when that feature has been reached.\n#\n# This change marks five more rows as reached by writing a lowercase\n# \"x\" into their \"Atteint\" cell:\n#   - four rows whose \"Atteint\" cell was empty gain the \"x\" text;\n#   - \"Tour \u00e0 4 joueurs\" already had an uppercase \"X\", normalized to \"x\".\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$targetLabels = @(\n    \"2 joueurs en r\u00e9seau\",\n    \"Affichage du plateau pour les 2 joueurs r\u00e9seau\",\n    \"Tour \u00e0 4 joueurs\",\n    \"4 joueurs en r\u00e9seau\",\n    \"Affichage du plateau pour les 4 joueurs r\u00e9seau\"\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    # Cell.Range.Text carries a trailing cell-end mark (CR + BEL, chars\n    # 13/7) that isn't whitespace, so strip those explicitly before\n    # comparing against the plain row label.\n    $label = $t.Cell($r, 1).Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($targetLabels -contains $label) {\n        $t.Cell($r, 2).Range.Text = \"x\"\n    }\n}\n"}
